$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.628.22'
$ws.Range("E2").Value = '  +1.74%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.458.68'
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.60'
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.14'
$ws.Range("E6").Value = '  +3.49%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.481'
$ws.Range("E8").Value = '  +1.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.61'
$ws.Range("E9").Value = '  -0.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.124'
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.399'
$ws.Range("E11").Value = '  +3.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.048.21'
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.70'
$ws.Range("E13").Value = '  +4.50%  '
$ws.Range("E14").Value = '  +2.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.471.82'
$ws.Range("E15").Value = '  +2.09%  '
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.754.23'
$ws.Range("E17").Value = '  +1.89%  '
$ws.Range("E18").Value = '  +3.81%  '
$ws.Range("E19").Value = '  +5.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.21'
$ws.Range("E20").Value = '  +2.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '388.40'
$ws.Range("E21").Value = '  -0.69%  '
$ws.Range("E22").Value = '  +2.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '74.82'
$ws.Range("E23").Value = '  -0.27%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.605.15'
$ws.Range("E25").Value = '  +2.20%  '
$ws.Range("E26").Value = '  +1.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.179'
$ws.Range("E27").Value = '  -7.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.57'
$ws.Range("E28").Value = '  +4.21%  '
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("E30").Value = '  +1.25%  '
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("E33").Value = '  -0.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.75'
$ws.Range("E34").Value = '  +2.09%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.27'
$ws.Range("E35").Value = '  +5.30%  '
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.07'
$ws.Range("E36").Value = '  +2.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '170.66'
$ws.Range("E37").Value = '  +1.16%  '
$ws.Range("E38").Value = '  +6.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '31.36'
$ws.Range("E39").Value = '  +21.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.497.87'
$ws.Range("E40").Value = '  +2.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0768'
$ws.Range("E41").Value = '  +0.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.796'
$ws.Range("E42").Value = '  +2.12%  '
$ws.Range("E43").Value = '  +1.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.13'
$ws.Range("E44").Value = '  -0.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.70'
$ws.Range("E45").Value = '  +3.35%  '
$ws.Range("E46").Value = '  +1.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.591.33'
$ws.Range("E47").Value = '  +4.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.28'
$ws.Range("E48").Value = '  +1.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.76'
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  +0.11%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.17'
$ws.Range("E51").Value = '  +6.54%  '
